# MIRO_setup.pptx - slide 8, "TextBox 4" shape: split the single
# "ln -s ~/lib/mdk-170906 ~/mdk " run into several runs (the MDK version
# number changed from 170906 to 180509), per the commit
# "Added MiRo touch control python code".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# Bump the date-stamped version embedded in the command (170906 -> 180509).
# Every other character in the run stays the same, so no earlier/later
# character offsets shift around this edit.
$tr.Characters(20, 5).Text = "80509"

# Re-assigning each slice's own (unchanged) text forces PowerPoint to cut a
# fresh run at that boundary, matching the finer-grained run layout in the
# target deck while leaving the visible formatting (font/size) untouched.
$tr.Characters(3, 8).Text = "ln -s ~/"
$tr.Characters(11, 3).Text = "lib"
$tr.Characters(14, 6).Text = "/mdk-1"
$tr.Characters(25, 1).Text = " "
$tr.Characters(26, 6).Text = "~/mdk "
